{"js": "// Replace the literal date \"16 \u0441\u0435\u0440\u043f\u043d\u044f 2018\" embedded in the sentence\n// \"\u0412\u043a\u0430\u0437\u0430\u043d\u0435 \u0440\u0456\u0448\u0435\u043d\u043d\u044f \u0431\u0443\u043b\u043e \u043f\u043e\u0434\u0430\u043d\u043e \u0434\u0435\u0440\u0436\u0430\u0432\u043d\u043e\u043c\u0443 \u0440\u0435\u0454\u0441\u0442\u0440\u0430\u0442\u043e\u0440\u0443 16 \u0441\u0435\u0440\u043f\u043d\u044f 2018 \u0440\u043e\u043a\u0443\n// \u0434\u043b\u044f \u043f\u0440\u043e\u0432\u0435\u0434\u0435\u043d\u043d\u044f \u0440\u0435\u0454\u0441\u0442\u0440\u0430\u0446\u0456\u0439\u043d\u043e\u0457 \u0434\u0456\u0457 ...\" with the docxFiller placeholder\n// \"{5}\", splitting the original single run into three runs (identical\n// formatting) the way Word does when a sub-range of a run's text is\n// edited in place.\n\nconst body = context.document.body;\n\n// Find the date substring inside the sentence (matchCase keeps us from\n// accidentally matching something else).\nconst dateResults = body.search(\"16 \u0441\u0435\u0440\u043f\u043d\u044f 2018\", { matchCase: true, matchWholeWord: false });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length === 0) {\n  throw new Error('Could not find the text \"16 \u0441\u0435\u0440\u043f\u043d\u044f 2018\" to replace.');\n}\n\nconst dateRange = dateResults.items[0];\n\n// Replace the date text in place with the placeholder.\nconst placeholderRange = dateRange.insertText(\"{5}\", \"Replace\");\nawait context.sync();\n\n// Toggling a character formatting property on the freshly inserted range\n// and then restoring it forces Word to materialize it as its own run\n// (separate from the surrounding text runs), matching the 3-run split\n// produced by the original edit even though the visible formatting ends\n// up identical on all three runs.\nplaceholderRange.font.bold = true;\nawait context.sync();\n\nplaceholderRange.font.bold = false;\nawait context.sync();\n", "ps1": "# Replace the literal date \"16 \u0441\u0435\u0440\u043f\u043d\u044f 2018\" embedded in the sentence\n# \"\u0412\u043a\u0430\u0437\u0430\u043d\u0435 \u0440\u0456\u0448\u0435\u043d\u043d\u044f \u0431\u0443\u043b\u043e \u043f\u043e\u0434\u0430\u043d\u043e \u0434\u0435\u0440\u0436\u0430\u0432\u043d\u043e\u043c\u0443 \u0440\u0435\u0454\u0441\u0442\u0440\u0430\u0442\u043e\u0440\u0443 16 \u0441\u0435\u0440\u043f\u043d\u044f 2018 \u0440\u043e\u043a\u0443\n# \u0434\u043b\u044f \u043f\u0440\u043e\u0432\u0435\u0434\u0435\u043d\u043d\u044f \u0440\u0435\u0454\u0441\u0442\u0440\u0430\u0446\u0456\u0439\u043d\u043e\u0457 \u0434\u0456\u0457 ...\" with the docxFiller placeholder\n# \"{5}\", splitting the original single run into three runs (identical\n# formatting) the way Word does when a sub-range of a run's text is\n# edited in place.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"16 \u0441\u0435\u0440\u043f\u043d\u044f 2018\")\n\nif ($found) {\n    # Replace the date text in place with the placeholder.\n    $rng.Text = \"{5}\"\n\n    # Toggling a character formatting property on the freshly inserted\n    # range and then restoring it forces Word to materialize it as its\n    # own run (separate from the surrounding text runs), matching the\n    # 3-run split produced by the original edit even though the visible\n    # formatting ends up identical on all three runs.\n    $rng.Font.Bold = 1\n    $rng.Font.Bold = 0\n}\n"}
